$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width: F (col 6) from 21 to 19 characters ---
# COM ColumnWidth units are offset from the raw OOXML width by -0.8333333333333334
$ws.Columns.Item(6).ColumnWidth = 18.166666666666668

# --- Re-seat names for this week's ride assignments ---
$ws.Range("E2").Value = "Driver: Kaitlyn Kim"
$ws.Range("F2").Value = "Driver: Josh Paik"
$ws.Range("C3").Value = "Ben Kim"
$ws.Range("D3").Value = "Gabriel Ni"
$ws.Range("D4").Copy($ws.Range("E3"))
$ws.Range("E3").Value = "Lindsey Ro"
$ws.Range("L3").Copy($ws.Range("F3"))
$ws.Range("F3").Value = "Jane Yoo"
$ws.Range("G3").Value = "Kyle Hwang"
$ws.Range("M3").Value = "Joann Jung"
$ws.Range("C4").Value = "Daniel Kuo"
$ws.Range("D4").Copy($ws.Range("E4"))
$ws.Range("E4").Value = "Eugene Seo"
$ws.Range("L3").Copy($ws.Range("F4"))
$ws.Range("F4").Value = "Grace Sowon Park"
$ws.Range("G4").Value = "Cara Lee"
$ws.Range("K4").Value = "Gabriel Ni"
$ws.Range("L4").Value = "Grace Sowon Park"
$ws.Range("M4").Value = "Sam Ko"
$ws.Range("D5").Value = "Sam Ko"
$ws.Range("L3").Copy($ws.Range("E5"))
$ws.Range("E5").Value = "helena song🐟"
$ws.Range("B5").Copy($ws.Range("F5"))
$ws.Range("F5").Value = "Joanna Wei"
$ws.Range("G5").Value = "Sehyun Jung"
$ws.Range("K5").Value = "Zoe Li"
$ws.Range("L5").Value = "Jane Yoo"
$ws.Range("M5").Value = "Cara Lee"
$ws.Range("C6").Value = "Joann Jung"
$ws.Range("E6").Clear()
$ws.Range("B5").Copy($ws.Range("F6"))
$ws.Range("F6").Value = "Isabelle Li"
$ws.Range("G6").Value = "Joel Shim"
$ws.Range("L6").Value = "Isabelle Li"
$ws.Range("M6").Value = "Claire Doh"
$ws.Range("M12").Value = "Joanna Wei — No valid driver"
$ws.Range("M13").Value = "Sehyun Jung — No valid driver"
